$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-11 (years 2000-2009), shifting remaining rows up
$ws.Range("A2:F11").EntireRow.Delete()

# Add new row 13 for 2021年
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 1218326
$ws.Range("C13").Value = 324210
$ws.Range("D13").Value = 709456
$ws.Range("E13").Value = 8743661
$ws.Range("F13").Value = 2104041
